# Apply the weekly update described in the commit "Fruta / hortaliza, semanal".
# Row 75 (previously the 2021-01-13 "Primera" record) is refreshed to the new
# 2023-01-13 observation; the old row 75 values move down to a new row 77, a
# new "Segunda" quality record is inserted as row 76, and the former row 76
# (2022-08-18 record) moves down to row 78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75: update in place with the new 2023-01-13 "Primera" observation ---
$ws.Range("D75").Value = 44939
$ws.Range("J75").Value = 200
$ws.Range("K75").Value = 600
$ws.Range("L75").Value = 600
$ws.Range("M75").Value = 600
$ws.Range("N75").Value = "$/paquete 6 unidades"
$ws.Range("P75").Value = 100
$ws.Range("Q75").Value = 6

# --- Row 76: new "Segunda" quality observation for 2023-01-13 ---
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44939
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = 100112037
$ws.Range("G76").Value = "Cebollín"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Segunda"
$ws.Range("J76").Value = 200
$ws.Range("K76").Value = 500
$ws.Range("L76").Value = 500
$ws.Range("M76").Value = 500
$ws.Range("N76").Value = "$/paquete 6 unidades"
$ws.Range("O76").Value = "Provincia de Diguillín"
$ws.Range("P76").Value = 83
$ws.Range("Q76").Value = 6
$ws.Range("R76").Value = "Hortaliza"

# --- Row 77: former row 75 data (2021-01-13, "Primera") moved down ---
$ws.Range("A77").Value = 7
$ws.Range("B77").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C77").Value = "Ñuble"
$ws.Range("D77").Value = 44209
$ws.Range("E77").Value = 16
$ws.Range("F77").Value = 100112037
$ws.Range("G77").Value = "Cebollín"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 150
$ws.Range("K77").Value = 3500
$ws.Range("L77").Value = 4000
$ws.Range("M77").Value = 3767
$ws.Range("N77").Value = "$/paquete 2 kilos"
$ws.Range("O77").Value = "Provincia de Diguillín"
$ws.Range("P77").Value = 1884
$ws.Range("Q77").Value = 2
$ws.Range("R77").Value = "Hortaliza"

# --- Row 78: former row 76 data (2022-08-18, "Primera") moved down ---
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44791
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112037
$ws.Range("G78").Value = "Cebollín"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 120
$ws.Range("K78").Value = 8000
$ws.Range("L78").Value = 8500
$ws.Range("M78").Value = 8250
$ws.Range("N78").Value = "$/docena de atados"
$ws.Range("O78").Value = "Provincia de Diguillín"
$ws.Range("P78").Value = 2750
$ws.Range("Q78").Value = 3
$ws.Range("R78").Value = "Hortaliza"

# Match the formatted-date style used by the rest of column D (numFmtId 165)
$ws.Range("D76").NumberFormat = $ws.Range("D74").NumberFormat()
$ws.Range("D77").NumberFormat = $ws.Range("D74").NumberFormat()
$ws.Range("D78").NumberFormat = $ws.Range("D74").NumberFormat()

Write-Output "Done"
